$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column G (K column), rows 2-19
$values = @{
    2  = 0
    3  = 0
    4  = 3
    5  = 4
    6  = 2
    7  = 3
    8  = 1
    9  = 2
    10 = 4
    11 = 1
    12 = 3
    13 = 3
    14 = 2
    15 = 3
    16 = 2
    17 = 4
    18 = 2
    19 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
